$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.049.80"
$ws.Range("D3").Value = "1.828.66"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "311.56"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "0.4334"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.07310"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "0.8445"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "20.67"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "1.827.69"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'6.670"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "5.295"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "89.63"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "0.000008782"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "14.93"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "27.099.69"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").Value = "5.144"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "10.89"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "2.051.21"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "1.983"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").Value = "151.24"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").Value = "2.214"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").Value = "'18.30"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "5.244"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "117.18"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "0.08732"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "'1.180"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "0.7404"
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("D34").Value = "4.443"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "2.906"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "0.9999"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "0.05238"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "7.228"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "0.1702"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "0.5138"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "8.578"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "10.64"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "0.4764"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "1.946"
$ws.Range("E47").Value = "  +6.38%  "
$ws.Range("D48").Value = "'105.90"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "0.9997"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "1.668"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  -1.58%  "
